$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.164.17"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.848.99"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'235.30"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4712"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").Value = "'0.2890"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").Value = "'0.06512"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'21.56"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").Value = "'0.07944"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "'97.35"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "1.852.24"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "'5.074"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "'0.6729"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'266.59"
$ws.Range("E16").Value = "  -4.88%  "
$ws.Range("D17").Value = "30.133.81"
$ws.Range("D18").Value = "'13.55"
$ws.Range("E18").Value = "  +6.91%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'0.000007526"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("D21").Value = "2.097.83"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'5.217"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").Value = "'6.121"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").Value = "'166.10"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "'9.142"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "'18.79"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'1.919"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "'1.394"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").Value = "'0.09832"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("D31").Value = "'1.462"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "'4.259"
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").Value = "'3.988"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("D34").Value = "'0.04665"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").Value = "'1.114"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'0.6948"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "'2.711"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "'0.01857"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "'2.605"
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("D40").Value = "'6.306"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Value = "'73.04"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'1.922"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "'103.07"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "'0.4100"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").Value = "'939.88"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").Value = "'9.070"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "'6.942"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").Value = "'33.66"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").Value = "'0.05654"
$ws.Range("E51").Value = "  +0.40%  "
